# Generate Report for Handoff
# Updates the localization-status workbook: a new source-file GUID/hash
# replaces the old one, handoff/handback timestamps are refreshed, the
# "Latest Target File" / "Latest Handback File" values are cleared on the
# per-locale sheets (no handback received yet), and the corresponding
# column widths shrink to fit the now-empty columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "e52db018-6830-4abf-ba72-dd23dfc01521"
$newGuid = "191efabe-4ebd-4ca0-95d5-9be48652d436"
$newHash = "488bcab9e837d8978e09db0808131cc4533dc3af"

$overviewHyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/$oldGuid.md"
$zhcnHyperlinkAddr     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/$oldGuid.md"
$dedeHyperlinkAddr     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-21 03:03:02"

# Refresh the single hyperlink (File Name -> Path And Name) with the new
# display text, keeping the same target address.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewHyperlinkAddr, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 03:02:56"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Only the Source File Name (A2) hyperlink remains; the Latest Target
# File (I2) hyperlink is removed along with its value.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnHyperlinkAddr, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Style = "Normal"

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 03:03:02"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeHyperlinkAddr, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Style = "Normal"

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
